$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 40
$ws_ALC.Range("H40").Value = 5955.2856
$ws_ALC.Range("I40").Value = 11106.1
$ws_ALC.Range("J40").Value = 1272.7273
$ws_ALC.Range("K40").Value = 11106.1
$ws_ALC.Range("L40").Value = 1272.7273
$ws_ALC.Range("M40").Value = -10931.1
$ws_ALC.Range("N40").Value = -1622.7273

# ALC row 41
$ws_ALC.Range("H41").Value = 228.04546
$ws_ALC.Range("I41").Value = 122.1
$ws_ALC.Range("J41").Value = 316.33334
$ws_ALC.Range("K41").Value = 122.1
$ws_ALC.Range("L41").Value = 316.33334
$ws_ALC.Range("M41").Value = 317.9
$ws_ALC.Range("N41").Value = -1196.33334

# ALC row 48
$ws_ALC.Range("H48").Value = 7798.364
$ws_ALC.Range("I48").Value = 8798.666999999999
$ws_ALC.Range("J48").Value = 6598
$ws_ALC.Range("K48").Value = 26396.001
$ws_ALC.Range("L48").Value = 19794
$ws_ALC.Range("M48").Value = -26104.001
$ws_ALC.Range("N48").Value = -20378

# ALC row 56
$ws_ALC.Range("H56").Value = 7798.364
$ws_ALC.Range("I56").Value = 8798.666999999999
$ws_ALC.Range("J56").Value = 6598
$ws_ALC.Range("K56").Value = 26396.001
$ws_ALC.Range("L56").Value = 19794
$ws_ALC.Range("M56").Value = -25862.001
$ws_ALC.Range("N56").Value = -20862

# ALC row 98
$ws_ALC.Range("H98").Value = 2279.2144
$ws_ALC.Range("I98").Value = 2164.4546
$ws_ALC.Range("J98").Value = 2700
$ws_ALC.Range("K98").Value = 2164.4546
$ws_ALC.Range("L98").Value = 2700
$ws_ALC.Range("M98").Value = -666.4546
$ws_ALC.Range("N98").Value = -5696

# ALC row 122
$ws_ALC.Range("H122").Value = 2279.2144
$ws_ALC.Range("I122").Value = 2164.4546
$ws_ALC.Range("J122").Value = 2700
$ws_ALC.Range("K122").Value = 6493.3638
$ws_ALC.Range("L122").Value = 8100
$ws_ALC.Range("M122").Value = -4043.3638
$ws_ALC.Range("N122").Value = -13000

# ARM row 2
$ws_ARM.Range("H2").Value = 956.4103
$ws_ARM.Range("I2").Value = 862.5
$ws_ARM.Range("J2").Value = 1195.4546
$ws_ARM.Range("K2").Value = 862.5
$ws_ARM.Range("L2").Value = 1195.4546
$ws_ARM.Range("M2").Value = -749.5
$ws_ARM.Range("N2").Value = -1421.4546

# ARM row 32
$ws_ARM.Range("H32").Value = 7757.5483
$ws_ARM.Range("I32").Value = 3910.7964
$ws_ARM.Range("J32").Value = 33723.125
$ws_ARM.Range("K32").Value = 3910.7964
$ws_ARM.Range("L32").Value = 33723.125
$ws_ARM.Range("M32").Value = -3623.7964
$ws_ARM.Range("N32").Value = -34297.125

# ARM row 60
$ws_ARM.Range("H60").Value = 16333.333
$ws_ARM.Range("I60").Value = 12000
$ws_ARM.Range("J60").Value = 25000
$ws_ARM.Range("K60").Value = 12000
$ws_ARM.Range("L60").Value = 25000
$ws_ARM.Range("M60").Value = -11267
$ws_ARM.Range("N60").Value = -26466

# ARM row 74
$ws_ARM.Range("H74").Value = 26316976
$ws_ARM.Range("I74").Value = 33334258
$ws_ARM.Range("J74").Value = 2173.5
$ws_ARM.Range("K74").Value = 33334258
$ws_ARM.Range("L74").Value = 2173.5
$ws_ARM.Range("M74").Value = -33333384
$ws_ARM.Range("N74").Value = -3921.5

# ARM row 77
$ws_ARM.Range("H77").Value = 26316976
$ws_ARM.Range("I77").Value = 33334258
$ws_ARM.Range("J77").Value = 2173.5
$ws_ARM.Range("K77").Value = 166671290
$ws_ARM.Range("L77").Value = 10867.5
$ws_ARM.Range("M77").Value = -166666922
$ws_ARM.Range("N77").Value = -19603.5

# ARM row 116
$ws_ARM.Range("H116").Value = 956.4103
$ws_ARM.Range("I116").Value = 862.5
$ws_ARM.Range("J116").Value = 1195.4546
$ws_ARM.Range("K116").Value = 862.5
$ws_ARM.Range("L116").Value = 1195.4546
$ws_ARM.Range("M116").Value = 1431.5
$ws_ARM.Range("N116").Value = -5783.4546

# BSM row 3
$ws_BSM.Range("H3").Value = 956.4103
$ws_BSM.Range("I3").Value = 862.5
$ws_BSM.Range("J3").Value = 1195.4546
$ws_BSM.Range("K3").Value = 862.5
$ws_BSM.Range("L3").Value = 1195.4546
$ws_BSM.Range("M3").Value = -748.5
$ws_BSM.Range("N3").Value = -1423.4546

# BSM row 82
$ws_BSM.Range("H82").Value = 22228.572
$ws_BSM.Range("I82").Value = 0
$ws_BSM.Range("J82").Value = 22228.572
$ws_BSM.Range("K82").Value = 0
$ws_BSM.Range("L82").Value = 22228.572
$ws_BSM.Range("M82").ClearContents()
$ws_BSM.Range("N82").Value = -22994.572

# BSM row 85
$ws_BSM.Range("H85").Value = 22228.572
$ws_BSM.Range("I85").Value = 0
$ws_BSM.Range("J85").Value = 22228.572
$ws_BSM.Range("K85").Value = 0
$ws_BSM.Range("L85").Value = 22228.572
$ws_BSM.Range("M85").ClearContents()
$ws_BSM.Range("N85").Value = -24880.572

# BSM row 134
$ws_BSM.Range("H134").Value = 14706959
$ws_BSM.Range("I134").Value = 20834400
$ws_BSM.Range("J134").Value = 1100
$ws_BSM.Range("K134").Value = 62503200
$ws_BSM.Range("L134").Value = 3300
$ws_BSM.Range("M134").Value = -62500665
$ws_BSM.Range("N134").Value = -8370

# CRP row 31
$ws_CRP.Range("H31").Value = 14708994
$ws_CRP.Range("I31").Value = 21741156
$ws_CRP.Range("J31").Value = 5383.364
$ws_CRP.Range("K31").Value = 21741156
$ws_CRP.Range("L31").Value = 5383.364
$ws_CRP.Range("M31").Value = -21740861
$ws_CRP.Range("N31").Value = -5973.364

# CRP row 34
$ws_CRP.Range("H34").Value = 14708994
$ws_CRP.Range("I34").Value = 21741156
$ws_CRP.Range("J34").Value = 5383.364
$ws_CRP.Range("K34").Value = 21741156
$ws_CRP.Range("L34").Value = 5383.364
$ws_CRP.Range("M34").Value = -21740954
$ws_CRP.Range("N34").Value = -5787.364

# CRP row 54
$ws_CRP.Range("H54").Value = 12000
$ws_CRP.Range("J54").Value = 12000
$ws_CRP.Range("L54").Value = 12000
$ws_CRP.Range("N54").Value = -13316

# CRP row 122
$ws_CRP.Range("H122").Value = 2038
$ws_CRP.Range("I122").Value = 1500
$ws_CRP.Range("J122").Value = 2360.8
$ws_CRP.Range("K122").Value = 4500
$ws_CRP.Range("L122").Value = 7082.400000000001
$ws_CRP.Range("M122").Value = -2050
$ws_CRP.Range("N122").Value = -11982.4

# CRP row 127
$ws_CRP.Range("H127").Value = 49999.5
$ws_CRP.Range("J127").Value = 49999.5
$ws_CRP.Range("L127").Value = 49999.5
$ws_CRP.Range("N127").Value = -59919.5

# CRP row 134
$ws_CRP.Range("H134").Value = 2218.8333
$ws_CRP.Range("I134").Value = 1416.8572
$ws_CRP.Range("J134").Value = 3341.6
$ws_CRP.Range("K134").Value = 4250.571599999999
$ws_CRP.Range("L134").Value = 10024.8
$ws_CRP.Range("M134").Value = -1715.571599999999
$ws_CRP.Range("N134").Value = -15094.8

# GSM row 123
$ws_GSM.Range("H123").Value = 10636.5
$ws_GSM.Range("J123").Value = 10636.5
$ws_GSM.Range("L123").Value = 10636.5
$ws_GSM.Range("N123").Value = -15536.5

# WVR row 136
$ws_WVR.Range("H136").Value = 8072.706
$ws_WVR.Range("I136").Value = 2027.4445
$ws_WVR.Range("J136").Value = 14873.625
$ws_WVR.Range("K136").Value = 6082.333500000001
$ws_WVR.Range("L136").Value = 44620.875
$ws_WVR.Range("M136").Value = -3532.333500000001
$ws_WVR.Range("N136").Value = -49720.875
